$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").ClearContents()
$ws.Range("H3").ClearContents()

$ws.Range("E4").ClearContents()
$ws.Range("H4").ClearContents()

$ws.Range("E5").ClearContents()
$ws.Range("H5").ClearContents()

$ws.Range("E6").ClearContents()
$ws.Range("H6").ClearContents()

$ws.Range("E8").ClearContents()
$ws.Range("H8").ClearContents()

$ws.Range("E9").ClearContents()
$ws.Range("H9").ClearContents()

$ws.Range("E10").ClearContents()
$ws.Range("H10").ClearContents()

$ws.Range("E11").ClearContents()
$ws.Range("H11").ClearContents()

$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 909,5 TL"
$ws.Range("E13").ClearContents()
$ws.Range("H13").ClearContents()

$ws.Range("E14").ClearContents()
$ws.Range("H14").ClearContents()
